# DataTables generation based on DocumentUnderstanding OCR
# Adds jobIndexPath / cvIndexPath / jobFolderPath / cvFolderPath settings
# rows to the "Constants" sheet, pushing the existing
# replyMessageNoAttachment / replyMessageConfirmation / apiKey rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Make room: insert 4 blank rows starting at row 13. This shifts the
# existing rows 14 (replyMessageNoAttachment), 15 (replyMessageConfirmation)
# and 17 (apiKey) down to 18, 19 and 21 respectively - matching the target
# layout (and keeping the existing blank-row gaps).
$ws.Range("A13:A16").EntireRow.Insert()

# Write the new "Name" column cells first, in the order that matches the
# desired shared-string insertion order (jobIndexPath, cvIndexPath,
# jobFolderPath, cvFolderPath), then the "Value" column cells (Jobs folder,
# CVs folder, jobs index file, CVs index file).
$ws.Range("A15").Value = "jobIndexPath"
$ws.Range("A16").Value = "cvIndexPath"
$ws.Range("A13").Value = "jobFolderPath"
$ws.Range("A14").Value = "cvFolderPath"

$ws.Range("B13").Value = "Data\Index\Jobs"
$ws.Range("B14").Value = "Data\Index\CVs"
$ws.Range("B15").Value = "Data\Index\Jobs\jobs.xlsx"
$ws.Range("B16").Value = "Data\Index\CVs\CVs.xlsx"

# Widen column B so the longer paths fit (engine quantizes stored width to
# 1/6-character steps, so 95.27 is the closest input to the target
# 96.109375 serialized width), and refresh the active selection (this also
# clears the stale topLeftCell="A10" scroll position).
$ws.Columns.Item(2).ColumnWidth = 95.27

[void]$ws.Activate()
[void]$ws.Range("B25").Select()
